# Apply the "Automatic update of files" change to rows 12-17 of the Artfynd sheet.
# The edit re-sequences which observation record occupies which row (columns
# A, B, D, E, F, G, H, I, J, Q, R, AC), refreshing the Taxonsorteringsordning (B)
# values to new values from the source at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 <- (previously row 16's species data)
$ws.Range("A12").Value = 112128530
$ws.Range("B12").Value = 89936
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 5420
$ws.Range("F12").Value = "Grovticka"
$ws.Range("G12").Value = "Phaeolus schweinitzii"
$ws.Range("H12").Value = "(Fr.) Pat."
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("Q12").Value = 657144
$ws.Range("R12").Value = 6571278
$ws.Range("AC12").Value = "På högstubbe av tall."

# Row 13 <- (previously row 12's species data)
$ws.Range("A13").Value = 112128551
$ws.Range("B13").Value = 90800
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 4364
$ws.Range("F13").Value = "Dropptaggsvamp"
$ws.Range("G13").Value = "Hydnellum ferrugineum"
$ws.Range("H13").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("Q13").Value = 657162
$ws.Range("R13").Value = 6571271
$ws.Range("AC13").Value = ""

# Row 14 <- (previously row 17's species data)
$ws.Range("A14").Value = 112128708
$ws.Range("B14").Value = 90821
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 5964
$ws.Range("F14").Value = "Fjällig taggsvamp s.str."
$ws.Range("G14").Value = "Sarcodon imbricatus s.str."
$ws.Range("H14").Value = "(L.:Fr.) P.Karst."
# Antal ("1") is stored as text in the source data, not a number - force text format
# so it round-trips as a string rather than being coerced to numeric.
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "1"
$ws.Range("J14").Value = "fruktkroppar"
$ws.Range("Q14").Value = 657216
$ws.Range("R14").Value = 6571313
$ws.Range("AC14").Value = ""

# Row 15 <- (previously row 14's species data)
$ws.Range("A15").Value = 112128627
$ws.Range("B15").Value = 90221
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 3298
$ws.Range("F15").Value = "Trådticka"
$ws.Range("G15").Value = "Climacocystis borealis"
$ws.Range("H15").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I15").Value = ""
$ws.Range("J15").Value = ""
$ws.Range("Q15").Value = 657182
$ws.Range("R15").Value = 6571192
$ws.Range("AC15").Value = "På nedre delen av torrgran."

# Row 16 <- (previously row 13's species data)
$ws.Range("A16").Value = 112128573
$ws.Range("B16").Value = 90802
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 788
$ws.Range("F16").Value = "Gul taggsvamp"
$ws.Range("G16").Value = "Hydnellum geogenium"
$ws.Range("H16").Value = "(Fr.) Banker"
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("Q16").Value = 657134
$ws.Range("R16").Value = 6571219
$ws.Range("AC16").Value = ""

# Row 17 <- (previously row 15's species data)
$ws.Range("A17").Value = 112128498
$ws.Range("B17").Value = 90152
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 1339
$ws.Range("F17").Value = "Brandticka"
$ws.Range("G17").Value = "Pycnoporellus fulgens"
$ws.Range("H17").Value = "(Fr.) Donk"
$ws.Range("I17").Value = ""
$ws.Range("J17").Value = ""
$ws.Range("Q17").Value = 657134
$ws.Range("R17").Value = 6571271
$ws.Range("AC17").Value = "På granlåga. En del årsfärska dödade granar av granbarkborre. Gott om död ved i form av torrträd och lågor av gran."
